$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.702.91"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.402.81"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "407.99"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.19"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.613"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.716"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.133"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -8.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.11"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.931.82"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.98"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000205"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -6.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.26"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.392.90"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.16"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.84%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "61.667.59"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "482.38"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +23.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "88.94"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.20"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.09"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.85%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "33.23"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.06"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.79"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.71"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.70"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.69"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.165"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.73%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -6.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.80"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -6.81%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "56.20"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0479"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.29%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "149.58"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +6.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.32"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.133"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.315"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.89"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.04"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.51"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.12"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.33"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +17.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "16.11"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "21.87"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.143"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +9.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "113.02"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +16.31%  "
